{"js": "// Update the division-problem table: replace each \"before\" expression\n// text with its corresponding \"after\" expression text, preserving the\n// existing run formatting (font/size) of each cell.\nconst replacements = [\n  [\"137\u00f76=\", \"891\u00f73=\"],\n  [\"773\u00f79=\", \"227\u00f73=\"],\n  [\"774\u00f76=\", \"337\u00f76=\"],\n  [\"317\u00f72=\", \"909\u00f74=\"],\n  [\"526\u00f79=\", \"590\u00f78=\"],\n  [\"834\u00f77=\", \"916\u00f79=\"],\n  [\"710\u00f79=\", \"221\u00f72=\"],\n  [\"346\u00f79=\", \"258\u00f76=\"],\n  [\"878\u00f79=\", \"975\u00f72=\"],\n  [\"446\u00f77=\", \"587\u00f72=\"],\n  [\"810\u00f73=\", \"944\u00f77=\"],\n  [\"837\u00f78=\", \"175\u00f74=\"],\n  [\"245\u00f73=\", \"705\u00f75=\"],\n  [\"346\u00f74=\", \"750\u00f79=\"],\n  [\"464\u00f78=\", \"779\u00f75=\"],\n  [\"395\u00f79=\", \"171\u00f74=\"],\n  [\"834\u00f76=\", \"839\u00f79=\"],\n  [\"990\u00f79=\", \"156\u00f77=\"],\n  [\"542\u00f77=\", \"193\u00f77=\"],\n  [\"372\u00f78=\", \"730\u00f72=\"],\n  [\"416\u00f79=\", \"447\u00f74=\"],\n  [\"420\u00f74=\", \"985\u00f79=\"],\n  [\"681\u00f73=\", \"426\u00f78=\"],\n  [\"748\u00f74=\", \"808\u00f78=\"],\n  [\"318\u00f72=\", \"925\u00f75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division-problem table: replace each \"before\" expression\n# text with its corresponding \"after\" expression text, preserving the\n# existing run formatting (font/size) of each cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"137\u00f76=\", \"891\u00f73=\"),\n    @(\"773\u00f79=\", \"227\u00f73=\"),\n    @(\"774\u00f76=\", \"337\u00f76=\"),\n    @(\"317\u00f72=\", \"909\u00f74=\"),\n    @(\"526\u00f79=\", \"590\u00f78=\"),\n    @(\"834\u00f77=\", \"916\u00f79=\"),\n    @(\"710\u00f79=\", \"221\u00f72=\"),\n    @(\"346\u00f79=\", \"258\u00f76=\"),\n    @(\"878\u00f79=\", \"975\u00f72=\"),\n    @(\"446\u00f77=\", \"587\u00f72=\"),\n    @(\"810\u00f73=\", \"944\u00f77=\"),\n    @(\"837\u00f78=\", \"175\u00f74=\"),\n    @(\"245\u00f73=\", \"705\u00f75=\"),\n    @(\"346\u00f74=\", \"750\u00f79=\"),\n    @(\"464\u00f78=\", \"779\u00f75=\"),\n    @(\"395\u00f79=\", \"171\u00f74=\"),\n    @(\"834\u00f76=\", \"839\u00f79=\"),\n    @(\"990\u00f79=\", \"156\u00f77=\"),\n    @(\"542\u00f77=\", \"193\u00f77=\"),\n    @(\"372\u00f78=\", \"730\u00f72=\"),\n    @(\"416\u00f79=\", \"447\u00f74=\"),\n    @(\"420\u00f74=\", \"985\u00f79=\"),\n    @(\"681\u00f73=\", \"426\u00f78=\"),\n    @(\"748\u00f74=\", \"808\u00f78=\"),\n    @(\"318\u00f72=\", \"925\u00f75=\")\n)\n\nforeach ($pair in $replacements) {\n    $before = $pair[0]\n    $after = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($before, $true, $false, $false, $false, $false, $true, 1, $false, $after, 2)\n}\n"}
